# Updated symbol list on Thu Feb  9 07:23:47 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto rows
# that moved since the last snapshot. Values are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string cell format) instead of auto-converting the numeric-looking
# strings / percentages into real numbers; the Style reset afterwards clears
# the "number stored as text" quote-prefix formatting Excel applies so the
# cell keeps its original (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}


# Row 2
Set-TextValue $ws "D2" "321.81"
Set-TextValue $ws "E2" "-3.02%"

# Row 3
Set-TextValue $ws "D3" "43.11"
Set-TextValue $ws "E3" "-6.11%"

# Row 4
Set-TextValue $ws "D4" "5.201"
Set-TextValue $ws "E4" "-7.26%"

# Row 5
Set-TextValue $ws "D5" "0.08189"

# Row 6
Set-TextValue $ws "D6" "4.326"
Set-TextValue $ws "E6" "-3.00%"

# Row 7
Set-TextValue $ws "D7" "1.834"
Set-TextValue $ws "E7" "-10.27%"

# Row 8
Set-TextValue $ws "D8" "0.9363"
Set-TextValue $ws "E8" "-4.02%"

# Row 9
Set-TextValue $ws "D9" "0.1113"
Set-TextValue $ws "E9" "-3.97%"

# Row 10
Set-TextValue $ws "D10" "0.1865"
Set-TextValue $ws "E10" "-3.11%"

# Row 11
Set-TextValue $ws "D11" "0.09348"
Set-TextValue $ws "E11" "-5.83%"

# Row 12
Set-TextValue $ws "D12" "0.04623"
Set-TextValue $ws "E12" "-0.92%"

# Row 13
Set-TextValue $ws "D13" "7.408"
Set-TextValue $ws "E13" "-28.66%"

# Row 14
Set-TextValue $ws "D14" "0.1057"
Set-TextValue $ws "E14" "-0.25%"

# Row 15
Set-TextValue $ws "D15" "0.001291"
Set-TextValue $ws "E15" "0.10%"

# Row 16
Set-TextValue $ws "D16" "0.005783"
Set-TextValue $ws "E16" "-4.27%"

# Row 17
Set-TextValue $ws "D17" "3.356"

# Row 18
Set-TextValue $ws "E18" "-1.84%"

# Row 19
Set-TextValue $ws "D19" "0.3363"
Set-TextValue $ws "E19" "0.01%"

# Row 20
Set-TextValue $ws "E20" "-0.27%"

# Row 21
Set-TextValue $ws "D21" "0.2623"
Set-TextValue $ws "E21" "-1.09%"

# Row 22
Set-TextValue $ws "D22" "0.04166"
Set-TextValue $ws "E22" "-0.68%"

# Row 23
Set-TextValue $ws "D23" "0.001249"
Set-TextValue $ws "E23" "-4.80%"

# Row 24
Set-TextValue $ws "D24" "0.004299"
Set-TextValue $ws "E24" "-6.95%"

# Row 25
Set-TextValue $ws "E25" "-15.62%"

# Row 26
Set-TextValue $ws "D26" "0.0002980"
Set-TextValue $ws "E26" "-20.48%"

# Row 38
Set-TextValue $ws "D38" "0.02720"
Set-TextValue $ws "E38" "-1.73%"

# Row 39
Set-TextValue $ws "D39" "0.05545"
Set-TextValue $ws "E39" "-4.45%"

# Row 40
Set-TextValue $ws "D40" "0.007975"
Set-TextValue $ws "E40" "3.28%"

# Row 41
Set-TextValue $ws "D41" "0.1393"
Set-TextValue $ws "E41" "-2.98%"

# Row 42
Set-TextValue $ws "D42" "0.006538"
Set-TextValue $ws "E42" "-10.19%"

# Row 43
Set-TextValue $ws "D43" "0.002090"
Set-TextValue $ws "E43" "3.80%"

# Row 44
Set-TextValue $ws "D44" "0.007459"
Set-TextValue $ws "E44" "-7.76%"

# Row 45
Set-TextValue $ws "D45" "0.3205"
Set-TextValue $ws "E45" "-5.82%"

# Row 46
Set-TextValue $ws "D46" "0.00006957"
Set-TextValue $ws "E46" "-4.77%"

# Row 47
Set-TextValue $ws "D47" "0.00000000749"
Set-TextValue $ws "E47" "-0.26%"

# Row 48
Set-TextValue $ws "D48" "0.003461"
Set-TextValue $ws "E48" "-1.00%"

# Row 49
Set-TextValue $ws "D49" "0.003531"
Set-TextValue $ws "E49" "0.75%"

# Row 50
Set-TextValue $ws "D50" "0.00002098"
Set-TextValue $ws "E50" "-0.26%"

# Row 51
Set-TextValue $ws "E51" "-0.26%"
